$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7: B7 becomes a numeric value instead of text
$ws.Range("B7").Value = 123456

# New rows 8-12
$ws.Range("A8").Value = "mark ruffalo"
$ws.Range("B8").Value = 40603
$ws.Range("C8").Value = $true
$ws.Range("D8").Value = 0

$ws.Range("A9").Value = "Elon Musk"
$ws.Range("B9").Value = "twitter"
$ws.Range("C9").Value = $true
$ws.Range("D9").Value = 0

$ws.Range("A10").Value = "Chris Evans"
$ws.Range("B10").Value = "hola"
$ws.Range("C10").Value = $true
$ws.Range("D10").Value = 0

$ws.Range("A11").Value = "vinicius"
$ws.Range("B11").NumberFormat = "@"
$ws.Range("B11").Value = "123456"
$ws.Range("B11").Style = "Normal"
$ws.Range("C11").Value = $true
$ws.Range("D11").Value = 0

$ws.Range("A12").Value = "Messi10"
$ws.Range("B12").NumberFormat = "@"
$ws.Range("B12").Value = "123456"
$ws.Range("B12").Style = "Normal"
$ws.Range("C12").Value = $true
$ws.Range("D12").Value = 0
